# Helper: write a value that Excel would otherwise auto-coerce to a number
# (e.g. "18") while still keeping the cell's plain (unstyled) text type, so
# the stored cell matches a user who typed a numeric-looking string into an
# already-text cell.
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "ID" -> "user_id"
$ws.Range("A1").Value = "user_id"

# Row 2: this visitor actually completed registration - fill in their real
# user id and profile answers (previously placeholder "g" values).
$ws.Range("A2").Value = 549874715
$ws.Range("C2").Value = "Марина боровик"
Set-TextValue $ws.Range("D2") "18"
$ws.Range("E2").Value = "Бентли"
$ws.Range("F2").Value = "Вроцлав"
Set-TextValue $ws.Range("G2") "2"

# Row 3: previously only had the id stamped (row created but unanswered) -
# fill in the rest of this registration.
$ws.Range("B3").Value = "Водитель"
$ws.Range("C3").Value = "Illia"
Set-TextValue $ws.Range("D3") "18"
$ws.Range("E3").Value = "ff"
$ws.Range("F3").Value = "fg"
$ws.Range("G3").Value = "ggg"

# Drop the trailing unused placeholder "Брокер" rows (old rows 4-7).
$ws.Rows.Item(4).EntireRow.Delete()
$ws.Rows.Item(4).EntireRow.Delete()
$ws.Rows.Item(4).EntireRow.Delete()
$ws.Rows.Item(4).EntireRow.Delete()

# Restore the saved selection state over the completed row.
$ws.Range("A3:G3").Select()
